# BountyPop review — move the "meta description" blurb from the top of the
# document down to the bottom (as a labelled "Play BountyPop Slot for Free -
# Review and Features" / body-text pair), replacing the old AI image-prompt
# paragraph that used to live at the very end of the document.
#
# Concretely, per the diff:
#   1. Delete the paragraph right after the H1 title that reads
#        "Meta description" (bold) + ": Read our BountyPop slot review ..."
#   2. Right before the final paragraph (the italic AI image-generation
#      prompt), insert a new bold paragraph reading
#        "Play BountyPop Slot for Free - Review and Features"
#   3. Replace the text of that final italic paragraph with
#        "Read our BountyPop slot review and find out the pros and cons to
#         play this high-volatility game for free."
#      (keeping its existing italic run formatting).

$d = $word.ActiveDocument

$oldMetaBody = "Create a cartoon-style feature image for BountyPop that features a happy Maya warrior with glasses. The image should be vibrant and eye-catching, featuring the Maya warrior surrounded by explosive gems and treasure chests, highlighting the adventurous pirate theme of the game. The image should encourage players to embrace the spirit of adventure and excitement while playing the game. The Maya warrior should be depicted as confident and carefree, perfectly embodying the attitude of players who are enjoying the game."
$newMetaBody = "Read our BountyPop slot review and find out the pros and cons to play this high-volatility game for free."
$newHeading  = "Play BountyPop Slot for Free - Review and Features"

# --- Step 1: replace the text of the final (italic) paragraph, keeping its formatting ---
$d.Content.Find.Execute($oldMetaBody, $true, $false, $false, $false, $false,
                         $true, 1, $false, $newMetaBody, 2)

# --- Step 2: build the new bold heading paragraph in a plain-formatted spot of the
#             document (so it doesn't inherit stray bold/italic direct formatting),
#             then relocate it (as FormattedText, mark included) to just before the
#             final paragraph. This keeps the run-level "bold" formatting intact
#             without leaking paragraph-mark formatting or neighbouring styles. ---

# Find a paragraph with plain run formatting to use as a scratch anchor: the one
# right before the "What we like" section heading.
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Some similar slots to BountyPop include Wild Pirates, Pirate Gold Deluxe, and Octopus Treasure.") {
        $anchorIndex = $i
        break
    }
}

$anchorRange = $d.Paragraphs.Item($anchorIndex).Range
$anchorRange.Collapse(0)
$anchorRange.InsertParagraphAfter()
$scratchIndex = $anchorIndex + 1
$scratchRange = $d.Paragraphs.Item($scratchIndex).Range

$bodyRange = $d.Range($scratchRange.Start, $scratchRange.End - 1)
$bodyRange.Text = $newHeading
$boldRange = $d.Range($scratchRange.Start, $scratchRange.Start + $newHeading.Length)
$boldRange.Font.Bold = $true

$scratchFormatted = $d.Paragraphs.Item($scratchIndex).Range.FormattedText

$lastIndex = $d.Paragraphs.Count
$lastRange = $d.Paragraphs.Item($lastIndex).Range
$destination = $d.Range($lastRange.Start, $lastRange.Start)
$destination.FormattedText = $scratchFormatted

# remove the scratch paragraph now that its formatted copy has been placed
$d.Paragraphs.Item($scratchIndex).Range.Delete()

# --- Step 3: delete the old "Meta description" paragraph near the top ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext.StartsWith("Meta description")) {
        $d.Paragraphs.Item($i).Range.Delete()
        break
    }
}
